# Update generated Excel file — add two new DBC message blocks
# (dir_actuator_feedback @ rows 12-18, dir_act_possition_loop @ rows 20-22)
# to the "autonomous_t26" sheet, and widen a handful of columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column width tweaks (stored OOXML width = ColumnWidth + 5/6, so we
#    subtract that offset to land exactly on the target stored widths).
# ---------------------------------------------------------------------
$offset = 5 / 6
$colWidths = @{ 1 = 33; 2 = 12; 3 = 19; 8 = 8; 9 = 7; 11 = 130 }
foreach ($col in $colWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col] - $offset
}

# ---------------------------------------------------------------------
# Helper data: header row (shared by every message block) and a couple
# of style helpers mirroring the sheet's existing look (message-title
# row / column-header row / plain data row).
# ---------------------------------------------------------------------
$headers = @("Signal Name", "Start Bit", "Length (bits)", "Byte Order", "Signed", "Factor", "Offset", "Min", "Max", "Unit", "Choices")

function Set-TitleRowStyle($range) {
    $range.Font.Bold = $true
    $range.Interior.Color = 15652797
    $range.Borders.LineStyle = 1
    $range.Borders.Weight = 2
}

function Set-HeaderRowStyle($range) {
    $range.Font.Bold = $true
    $range.Interior.Color = 6740479
    $range.Borders.LineStyle = 1
    $range.Borders.Weight = 2
}

function Set-DataRowStyle($range) {
    $range.Borders.LineStyle = 1
    $range.Borders.Weight = 2
}

# ---------------------------------------------------------------------
# 2. Message block: dir_actuator_feedback (rows 12-18)
# ---------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value = "Message: dir_actuator_feedback"
$ws.Cells.Item(12, 2).Value = "ID: 0x2968"
$ws.Cells.Item(12, 3).Value = "Sender(s): AK10_9"
Set-TitleRowStyle $ws.Range("A12:C12")

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(13, $i + 1).Value = $headers[$i]
}
Set-HeaderRowStyle $ws.Range("A13:K13")

$signals1 = @(
    @("Position", 0, 16, "Intel", $false, 10, 0, -32000, 32000, "º", $null),
    @("Speed", 16, 16, "Intel", $false, 10, 0, -32000, 32000, "RPM", $null),
    @("Current", 32, 16, "Intel", $false, 100, 0, $null, $null, "A", $null),
    @("Temperature", 48, 8, "Intel", $false, 1, 0, -20, 127, "º", $null),
    @("Error_codes", 56, 8, "Intel", $false, 1, 0, 0, 7, $null, "0=No Fault, 1=Motor Over-temperature, 2=Over-current, 3=Over-voltage, 4=Encode Fault, 5=Mosfet Over-temperature, 6=Motor Lock-up")
)

$row = 14
foreach ($sig in $signals1) {
    $ws.Cells.Item($row, 1).Value = $sig[0]
    $ws.Cells.Item($row, 2).Value = $sig[1]
    $ws.Cells.Item($row, 3).Value = $sig[2]
    $ws.Cells.Item($row, 4).Value = $sig[3]
    $ws.Cells.Item($row, 5).Value = $sig[4]
    $ws.Cells.Item($row, 6).Value = $sig[5]
    $ws.Cells.Item($row, 7).Value = $sig[6]
    if ($null -ne $sig[7]) { $ws.Cells.Item($row, 8).Value = $sig[7] }
    if ($null -ne $sig[8]) { $ws.Cells.Item($row, 9).Value = $sig[8] }
    if ($null -ne $sig[9]) { $ws.Cells.Item($row, 10).Value = $sig[9] }
    if ($null -ne $sig[10]) { $ws.Cells.Item($row, 11).Value = $sig[10] }
    Set-DataRowStyle $ws.Range("A" + $row + ":K" + $row)
    $row++
}

# ---------------------------------------------------------------------
# 3. Message block: dir_act_possition_loop (rows 20-22)
# ---------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value = "Message: dir_act_possition_loop"
$ws.Cells.Item(20, 2).Value = "ID: 0x468"
$ws.Cells.Item(20, 3).Value = "Sender(s): JETSON"
Set-TitleRowStyle $ws.Range("A20:C20")

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(21, $i + 1).Value = $headers[$i]
}
Set-HeaderRowStyle $ws.Range("A21:K21")

$ws.Cells.Item(22, 1).Value = "RPM"
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(22, 3).Value = 32
$ws.Cells.Item(22, 4).Value = "Intel"
$ws.Cells.Item(22, 5).Value = $true
$ws.Cells.Item(22, 6).Value = 10000
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 10).Value = "RPM"
Set-DataRowStyle $ws.Range("A22:K22")
